$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.867.32'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '1.809.95'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''309.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = '''0.4645'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.65%  '
$ws.Range('D8').Value = '''0.3694'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.43%  '
$ws.Range('D9').Value = '''0.07352'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').Value = '''0.8722'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = '  -1.60%  '
$ws.Range('D12').Value = '1.832.42'
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').Value = '''5.357'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('D14').Value = '''6.515'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.60%  '
$ws.Range('D15').Value = '''0.07048'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').Value = '''91.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.95%  '
$ws.Range('D17').Value = '''1.002'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '''0.000008705'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('D21').Value = '26.893.06'
$ws.Range('D22').Value = '''5.326'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = '''10.53'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.20%  '
$ws.Range('D24').Value = '2.069.11'
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('D25').Value = '''1.902'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.97%  '
$ws.Range('D26').Value = '''151.60'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').Value = '''2.143'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.86%  '
$ws.Range('D29').Value = '''5.311'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('D30').Value = '''115.88'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.14%  '
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('D32').Value = '''0.7557'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.08%  '
$ws.Range('E33').Value = '  -3.42%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''4.459'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.77%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '''2.918'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = '''1.098'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '''2.430'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.89%  '
$ws.Range('D41').Value = '''2.925'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').Value = '''0.5312'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('D43').Value = '''7.170'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('D44').Value = '''0.1664'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.33%  '
$ws.Range('D45').Value = '''8.447'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('D46').Value = '''0.4943'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.94%  '
$ws.Range('D47').Value = '''10.29'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.32%  '
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('D50').Value = '''103.11'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.60%  '
$ws.Range('E51').Value = '  -1.64%  '
